# Insert a new weekly record at row 30 ("Fruta / hortaliza, semanal").
# All existing rows from 30 down to 68 shift down by one (to 31..69),
# and the newly inserted row 30 receives the new raspberry price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(30).Insert()

$ws.Range("A30").Value = 9
$ws.Range("B30").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C30").Value = "Metropolitana"
$ws.Range("D30").Value = 44579
$ws.Range("E30").Value = 13
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100101
$ws.Range("H30").Value = "Berries"
$ws.Range("I30").Value = 100101004
$ws.Range("J30").Value = "Frambuesa"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 550
$ws.Range("N30").Value = 7500
$ws.Range("O30").Value = 8000
$ws.Range("P30").Value = 7809
$ws.Range("Q30").Value = "$/bandeja 2 kilos"
$ws.Range("R30").Value = "Provincia de Linares"
$ws.Range("S30").Value = 3904
$ws.Range("T30").Value = 2
